$d = $word.ActiveDocument

# --- Locate the bullet-list block to collapse -----------------------------
# The "Contenu pédagogique (synthèse)" section currently holds ~10 separate
# bulleted <w:p> paragraphs (one objective each). They all need to collapse
# into the single first paragraph of the block, whose text becomes the
# {{CONTENU_PEDAGOGIQUE}} merge field (paragraph formatting / numbering of
# that first paragraph must be preserved).
$startIndex = -1
$endIndex = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    $t = $p.Range.Text
    if ($startIndex -eq -1 -and $t -like "*cartographier les concepts*") {
        $startIndex = $i
    }
    if ($t -like "*la compl*tude.*") {
        $endIndex = $i
    }
}

if ($startIndex -eq -1 -or $endIndex -eq -1 -or $endIndex -le $startIndex) {
    Write-Host "Could not locate target block (start=$startIndex end=$endIndex)"
} else {
    # --- Remove every paragraph after the first one in the block ----------
    # (their bullets/text are folded away; only the first paragraph mark,
    # carrying the w:numPr bullet formatting, survives).
    $afterFirst = $d.Paragraphs($startIndex + 1)
    $lastOfBlock = $d.Paragraphs($endIndex)
    $killRange = $d.Range($afterFirst.Range.Start, $lastOfBlock.Range.End)
    $killRange.Delete()

    # --- Replace the remaining paragraph's own text with the merge field --
    $first = $d.Paragraphs($startIndex)
    $firstRange = $first.Range
    $firstRange.End = $firstRange.End - 1
    $insertStart = $firstRange.Start
    $firstRange.Text = ""

    $newRun = $d.Range($insertStart, $insertStart)
    $newRun.Text = "{{CONTENU_PEDAGOGIQUE}}"

    # --- Restyle: plain Calibri run (text was wiped above, so the old
    # bold/italic/font formatting isn't carried over into the new run) -----
    $first2 = $d.Paragraphs($startIndex)
    $finalRange = $first2.Range
    $finalRange.End = $finalRange.End - 1
    $finalRange.Font.Name = "Calibri"

    Write-Host "Collapsed paragraphs $startIndex..$endIndex into {{CONTENU_PEDAGOGIQUE}}"
}
